$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add repeated header row at row 42 (A42:E42) ---
$ws.Range("A42").Value = "name"
$ws.Range("B42").Value = "rh"
$ws.Range("C42").Value = "t"
$ws.Range("D42").Value = "lat"
$ws.Range("E42").Value = "lng"

# --- Fill in the B/C (rh/t Ubidots id) columns for the existing rows 43-46 ---
$ws.Range("B43").Value = "65d6457c7a715d000bf94dc0"
$ws.Range("C43").Value = "65d6457d7a715d000c7d068c"

$ws.Range("B44").Value = "65d644066eb306000dee50f9"
$ws.Range("C44").Value = "65d644067a7226000bcdb493"

$ws.Range("B45").Value = "65c3c37648bb6b000e4e9979"
$ws.Range("C45").Value = "65c3c377e858cb000eb4367a"

$ws.Range("B46").Value = "6644f0dd573ffb000ce2b97f"
$ws.Range("C46").Value = "6644f0ddb921b4000b90b85d"

# --- Update the "Date Delivered to Xin" values for rows 43-46 (7/12 -> 7/19) ---
$ws.Range("F43").Value = 45492
$ws.Range("F44").Value = 45492
$ws.Range("F45").Value = 45492
$ws.Range("F46").Value = 45492

# --- Append new sensor rows 47-49 ---
$ws.Range("A47").Value = "WS26-5KL"
$ws.Range("B47").Value = "65d644ba5ee5f8000c3ec157"
$ws.Range("C47").Value = "65d644bb7a7226000c15eb7f"
$ws.Range("D47").Value = 40.130597000000002
$ws.Range("E47").Value = -105.031846
$ws.Range("F47").Value = 45492
$ws.Range("G47").Value = "Ella Stankiewicz "
$ws.Range("H47").Value = "CO"

$ws.Range("A48").Value = "Drab-nb-rt-16"
$ws.Range("B48").Value = "N/A"
$ws.Range("C48").Value = "N/A"
$ws.Range("D48").Value = 44.896858999999999
$ws.Range("E48").Value = -108.561176
$ws.Range("F48").Value = 45492
$ws.Range("G48").Value = "Mark Bjornestad "
$ws.Range("H48").Value = "WY"

$ws.Range("A49").Value = "Drab-nb-rt-15"
$ws.Range("B49").Value = "N/A"
$ws.Range("C49").Value = "N/A"
$ws.Range("D49").Value = 44.204608999999998
$ws.Range("E49").Value = -107.92478300000001
$ws.Range("F49").Value = 45492
$ws.Range("G49").Value = "Mark Bjornestad "
$ws.Range("H49").Value = "WY"

# --- Match formatting of the surrounding data rows/cells ---
# Row 42 (repeated header): center aligned, General format (same as row 1 header)
$ws.Range("A42:E42").HorizontalAlignment = -4108

# Row 47: A/D/E/G/H centered General; F is a centered date (d-mmm-yy); B/C left unstyled
$ws.Range("A47").HorizontalAlignment = -4108
$ws.Range("D47:E47").HorizontalAlignment = -4108
$ws.Range("F47").HorizontalAlignment = -4108
$ws.Range("F47").NumberFormat = "d-mmm-yy"
$ws.Range("G47:H47").HorizontalAlignment = -4108

# Row 48: A/B/C/D/E/G/H centered General; F is a centered date (d-mmm-yy)
$ws.Range("A48:E48").HorizontalAlignment = -4108
$ws.Range("F48").HorizontalAlignment = -4108
$ws.Range("F48").NumberFormat = "d-mmm-yy"
$ws.Range("G48:H48").HorizontalAlignment = -4108

# Row 49: A/B/C/D/E/G/H centered General; F is a centered date (d-mmm-yy)
$ws.Range("A49:E49").HorizontalAlignment = -4108
$ws.Range("F49").HorizontalAlignment = -4108
$ws.Range("F49").NumberFormat = "d-mmm-yy"
$ws.Range("G49:H49").HorizontalAlignment = -4108

# Row 47/48/49 row height like the source (15pt custom height)
$ws.Range("47:49").RowHeight = 15
